# Insert a new weekly price record for "Jengibre" (Vega Modelo de Temuco)
# at the top of the data block (row 55), pushing all existing records
# (rows 55-167) down by one row (to rows 56-168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows down by inserting a new row at position 55.
$ws.Rows("55:55").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 44708
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = 100114007
$ws.Range("G55").Value = "Jengibre"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 30
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = 20000
$ws.Range("N55").Value = "`$/caja 13 kilos"
$ws.Range("O55").Value = "Perú"
$ws.Range("P55").Value = 1538
$ws.Range("Q55").Value = 13
$ws.Range("R55").Value = "Hortaliza"
